# Fruta / hortaliza, semanal
# Insert 3 new weekly records at the top of the "Vega Modelo de Temuco -
# Zapallo italiano" data block (rows 242-267), pushing the existing rows
# down to 245-270, then fill in the 3 new rows with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 242:267 down by 3 rows (new blank rows appear at 242:244)
$ws.Rows("242:244").Insert()

# Common (unchanged) values shared by every row in this data block
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 100112032
$categoria = "Zapallo italiano"
$variedad  = "Sin especificar"
$calidad   = "Primera"
$unidad    = "`$/caja 60 unidades"
$kgUnid    = 60
$clasif    = "Hortaliza"

function Set-Fila($fila, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg) {
    $ws.Cells.Item($fila, 1).Value  = $mercadoId
    $ws.Cells.Item($fila, 2).Value  = $mercado
    $ws.Cells.Item($fila, 3).Value  = $region
    $ws.Cells.Item($fila, 4).Value  = $fecha
    $ws.Cells.Item($fila, 5).Value  = $codreg
    $ws.Cells.Item($fila, 6).Value  = $catId
    $ws.Cells.Item($fila, 7).Value  = $categoria
    $ws.Cells.Item($fila, 8).Value  = $variedad
    $ws.Cells.Item($fila, 9).Value  = $calidad
    $ws.Cells.Item($fila, 10).Value = $volumen
    $ws.Cells.Item($fila, 11).Value = $precioMin
    $ws.Cells.Item($fila, 12).Value = $precioMax
    $ws.Cells.Item($fila, 13).Value = $precioProm
    $ws.Cells.Item($fila, 14).Value = $unidad
    $ws.Cells.Item($fila, 15).Value = $origen
    $ws.Cells.Item($fila, 16).Value = $precioKg
    $ws.Cells.Item($fila, 17).Value = $kgUnid
    $ws.Cells.Item($fila, 18).Value = $clasif
}

Set-Fila 242 44491 50 16000 16000 16000 "Limache" 267
Set-Fila 243 44491 80 16000 16000 16000 "Región de O'Higgins" 267
Set-Fila 244 44491 50 16000 16000 16000 "Región del Maule" 267
